$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    3 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    4 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    5 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    6 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 17.65757632934944)
    7 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B=2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
